$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (shifts IsLuring/Lokasi/Keterangan/FilePath one column right)
$ws.Columns.Item(5).Insert()

# The data that used to live in column E (IsLuring) lost its original number-format
# styling once it landed in column F - reset it back to the plain "Normal" style.
$ws.Range("F2:F4").Style = "Normal"

# New "Dosen" (lecturer) column header + per-row values (now supports multiple values e.g. "1,2")
$ws.Range("E1").Value = "Dosen"
$ws.Range("E2").Value = "1"
$ws.Range("E3").Value = "2"
$ws.Range("E4").Value = "1,2"

# Update the active selection to match the authored workbook
$ws.Range("E2").Select()
